$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data refresh drops "São Paulo" (old row 4) and "Mato
# Grosso" (old row 7), and introduces a new "Acre" row. Use real row
# delete/insert operations (shift cells) so that the untouched cells
# (notably the blank "Posição" cells for Brasil/Nordeste) keep flowing
# down with their original rows instead of being recreated from scratch.

# Remove "São Paulo" (row 4).
$ws.Rows(4).Delete()

# "Mato Grosso" has shifted up to row 6 now; remove it too.
$ws.Rows(6).Delete()

# Make room for the new "Acre" row ahead of "Tocantins" (now row 5).
$ws.Rows(5).Insert()

# Sheet now has 9 data+header rows (A1:D9). Rewrite the values/labels for
# 2024 vs 2013 comparison, re-ranked by the new "Preço médio" change.

# Row 2: Santa Catarina (rank 1, unchanged rank/position)
$ws.Cells.Item(2, 1).Value = "Santa Catarina"
$ws.Cells.Item(2, 2).Value = 36.18619365798551
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = "2024-2013"

# Row 3: Alagoas (rank 2)
$ws.Cells.Item(3, 1).Value = "Alagoas"
$ws.Cells.Item(3, 2).Value = 14.06202497955705
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = "2024-2013"

# Row 4: Sergipe (rank 3)
$ws.Cells.Item(4, 1).Value = "Sergipe"
$ws.Cells.Item(4, 2).Value = 9.903224112947802
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = "2024-2013"

# Row 5: Acre (rank 4, new row)
$ws.Cells.Item(5, 1).Value = "Acre"
$ws.Cells.Item(5, 2).Value = 8.673905343232624
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = "2024-2013"

# Row 6: Tocantins (rank 5)
$ws.Cells.Item(6, 1).Value = "Tocantins"
$ws.Cells.Item(6, 2).Value = 7.880778098376318
$ws.Cells.Item(6, 3).Value = 5
$ws.Cells.Item(6, 4).Value = "2024-2013"

# Row 7: Bahia (rank 6)
$ws.Cells.Item(7, 1).Value = "Bahia"
$ws.Cells.Item(7, 2).Value = 3.288256738162775
$ws.Cells.Item(7, 3).Value = 6
$ws.Cells.Item(7, 4).Value = "2024-2013"

# Row 8: Brasil (no ranking - the "Posição" cell stays blank, as before)
$ws.Cells.Item(8, 1).Value = "Brasil"
$ws.Cells.Item(8, 2).Value = -11.85100096624091
$ws.Cells.Item(8, 4).Value = "2024-2013"

# Row 9: Nordeste (no ranking - the "Posição" cell stays blank, as before)
$ws.Cells.Item(9, 1).Value = "Nordeste"
$ws.Cells.Item(9, 2).Value = -16.72189676051594
$ws.Cells.Item(9, 4).Value = "2024-2013"
